$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 348.17648
$ws.Cells.Item(28, 9).Value = 279.22223
$ws.Cells.Item(28, 10).Value = 425.75
$ws.Cells.Item(28, 11).Value = 279.22223
$ws.Cells.Item(28, 12).Value = 425.75
$ws.Cells.Item(28, 13).Value = 205.77777
$ws.Cells.Item(28, 14).Value = -1395.75

$ws.Cells.Item(107, 8).Value = 441.47058
$ws.Cells.Item(107, 9).Value = 427.66666
$ws.Cells.Item(107, 10).Value = 457
$ws.Cells.Item(107, 11).Value = 427.66666
$ws.Cells.Item(107, 12).Value = 457
$ws.Cells.Item(107, 13).Value = 1492.33334
$ws.Cells.Item(107, 14).Value = -4297

$ws.Cells.Item(110, 8).Value = 54350
$ws.Cells.Item(110, 10).Value = 54350
$ws.Cells.Item(110, 12).Value = 54350
$ws.Cells.Item(110, 14).Value = -62530

$ws.Cells.Item(118, 8).Value = 1023.5625
$ws.Cells.Item(118, 9).Value = 0
$ws.Cells.Item(118, 10).Value = 1023.5625
$ws.Cells.Item(118, 11).Value = 0
$ws.Cells.Item(118, 12).Value = 3070.6875
$ws.Cells.Item(118, 14).Value = -6384.6875
$ws.Cells.Item(118, 13).ClearContents()

$ws.Cells.Item(130, 8).Value = 49778.5
$ws.Cells.Item(130, 10).Value = 49778.5
$ws.Cells.Item(130, 12).Value = 49778.5
$ws.Cells.Item(130, 14).Value = -59818.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 776.9818
$ws.Cells.Item(2, 9).Value = 652.4865
$ws.Cells.Item(2, 10).Value = 1032.8889
$ws.Cells.Item(2, 11).Value = 652.4865
$ws.Cells.Item(2, 12).Value = 1032.8889
$ws.Cells.Item(2, 13).Value = -539.4865
$ws.Cells.Item(2, 14).Value = -1258.8889

$ws.Cells.Item(63, 8).Value = 1670.7778
$ws.Cells.Item(63, 9).Value = 1675.6
$ws.Cells.Item(63, 10).Value = 1657
$ws.Cells.Item(63, 11).Value = 1675.6
$ws.Cells.Item(63, 12).Value = 1657
$ws.Cells.Item(63, 13).Value = -989.5999999999999
$ws.Cells.Item(63, 14).Value = -3029

$ws.Cells.Item(66, 8).Value = 1670.7778
$ws.Cells.Item(66, 9).Value = 1675.6
$ws.Cells.Item(66, 10).Value = 1657
$ws.Cells.Item(66, 11).Value = 8378
$ws.Cells.Item(66, 12).Value = 8285
$ws.Cells.Item(66, 13).Value = -4946
$ws.Cells.Item(66, 14).Value = -15149

$ws.Cells.Item(74, 8).Value = 9091966
$ws.Cells.Item(74, 9).Value = 10639396
$ws.Cells.Item(74, 11).Value = 10639396
$ws.Cells.Item(74, 13).Value = -10638522

$ws.Cells.Item(77, 8).Value = 9091966
$ws.Cells.Item(77, 9).Value = 10639396
$ws.Cells.Item(77, 11).Value = 53196980
$ws.Cells.Item(77, 13).Value = -53192612

$ws.Cells.Item(110, 8).Value = 2080
$ws.Cells.Item(110, 9).Value = 1044
$ws.Cells.Item(110, 10).Value = 3116
$ws.Cells.Item(110, 11).Value = 1044
$ws.Cells.Item(110, 12).Value = 3116
$ws.Cells.Item(110, 13).Value = 1001
$ws.Cells.Item(110, 14).Value = -7206

$ws.Cells.Item(116, 8).Value = 776.9818
$ws.Cells.Item(116, 9).Value = 652.4865
$ws.Cells.Item(116, 10).Value = 1032.8889
$ws.Cells.Item(116, 11).Value = 652.4865
$ws.Cells.Item(116, 12).Value = 1032.8889
$ws.Cells.Item(116, 13).Value = 1641.5135
$ws.Cells.Item(116, 14).Value = -5620.8889

$ws.Cells.Item(132, 8).Value = 166685840
$ws.Cells.Item(132, 9).Value = 333336000
$ws.Cells.Item(132, 10).Value = 35666.668
$ws.Cells.Item(132, 11).Value = 1000008000
$ws.Cells.Item(132, 12).Value = 107000.004
$ws.Cells.Item(132, 13).Value = -1000005470
$ws.Cells.Item(132, 14).Value = -112060.004

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 776.9818
$ws.Cells.Item(3, 9).Value = 652.4865
$ws.Cells.Item(3, 10).Value = 1032.8889
$ws.Cells.Item(3, 11).Value = 652.4865
$ws.Cells.Item(3, 12).Value = 1032.8889
$ws.Cells.Item(3, 13).Value = -538.4865
$ws.Cells.Item(3, 14).Value = -1260.8889

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(68, 8).Value = 18397.5
$ws.Cells.Item(68, 10).Value = 18397.5
$ws.Cells.Item(68, 12).Value = 18397.5
$ws.Cells.Item(68, 14).Value = -19895.5

$ws.Cells.Item(71, 8).Value = 18397.5
$ws.Cells.Item(71, 10).Value = 18397.5
$ws.Cells.Item(71, 12).Value = 55192.5
$ws.Cells.Item(71, 14).Value = -62680.5

$ws.Cells.Item(93, 8).Value = 7635.6665
$ws.Cells.Item(93, 9).Value = 7635.6665
$ws.Cells.Item(93, 11).Value = 7635.6665
$ws.Cells.Item(93, 13).Value = -5763.6665

$ws.Cells.Item(94, 8).Value = 2211.6924
$ws.Cells.Item(94, 10).Value = 2229.3333
$ws.Cells.Item(94, 12).Value = 2229.3333
$ws.Cells.Item(94, 14).Value = -3131.3333

$ws.Cells.Item(132, 8).Value = 19615222
$ws.Cells.Item(132, 9).Value = 1388.8182
$ws.Cells.Item(132, 10).Value = 55573916
$ws.Cells.Item(132, 11).Value = 4166.4546
$ws.Cells.Item(132, 12).Value = 166721748
$ws.Cells.Item(132, 13).Value = -1636.4546
$ws.Cells.Item(132, 14).Value = -166726808

$ws.Cells.Item(134, 8).Value = 951.76
$ws.Cells.Item(134, 9).Value = 898.2105
$ws.Cells.Item(134, 10).Value = 1121.3334
$ws.Cells.Item(134, 11).Value = 2694.6315
$ws.Cells.Item(134, 12).Value = 3364.0002
$ws.Cells.Item(134, 13).Value = -159.6315
$ws.Cells.Item(134, 14).Value = -8434.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(87, 8).Value = 1160
$ws.Cells.Item(87, 9).Value = 314
$ws.Cells.Item(87, 10).Value = 2006
$ws.Cells.Item(87, 11).Value = 942
$ws.Cells.Item(87, 12).Value = 6018
$ws.Cells.Item(87, 13).Value = 306
$ws.Cells.Item(87, 14).Value = -8514

$ws.Cells.Item(90, 8).Value = 1160
$ws.Cells.Item(90, 9).Value = 314
$ws.Cells.Item(90, 10).Value = 2006
$ws.Cells.Item(90, 11).Value = 2826
$ws.Cells.Item(90, 12).Value = 18054
$ws.Cells.Item(90, 13).Value = 3414
$ws.Cells.Item(90, 14).Value = -30534

$ws.Cells.Item(113, 8).Value = 4825051.5
$ws.Cells.Item(113, 9).Value = 20833734
$ws.Cells.Item(113, 10).Value = 2941676.8
$ws.Cells.Item(113, 11).Value = 62501202
$ws.Cells.Item(113, 12).Value = 8825030.399999999
$ws.Cells.Item(113, 13).Value = -62499032
$ws.Cells.Item(113, 14).Value = -8829370.399999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 1868.6842
$ws.Cells.Item(113, 9).Value = 1979.4
$ws.Cells.Item(113, 10).Value = 1745.6666
$ws.Cells.Item(113, 11).Value = 1979.4
$ws.Cells.Item(113, 12).Value = 1745.6666
$ws.Cells.Item(113, 13).Value = 190.5999999999999
$ws.Cells.Item(113, 14).Value = -6085.6666

$ws.Cells.Item(126, 8).Value = 4931.6665
$ws.Cells.Item(126, 9).Value = 7550
$ws.Cells.Item(126, 10).Value = 3622.5
$ws.Cells.Item(126, 11).Value = 22650
$ws.Cells.Item(126, 12).Value = 10867.5
$ws.Cells.Item(126, 13).Value = -20180
$ws.Cells.Item(126, 14).Value = -15807.5

$ws.Cells.Item(136, 8).Value = 11769.152
$ws.Cells.Item(136, 10).Value = 11769.152
$ws.Cells.Item(136, 12).Value = 35307.456
$ws.Cells.Item(136, 14).Value = -40407.456

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 7814050
$ws.Cells.Item(40, 9).Value = 1608.08
$ws.Cells.Item(40, 10).Value = 35715628
$ws.Cells.Item(40, 11).Value = 1608.08
$ws.Cells.Item(40, 12).Value = 35715628
$ws.Cells.Item(40, 13).Value = -1472.08
$ws.Cells.Item(40, 14).Value = -35715900

$ws.Cells.Item(70, 8).Value = 29800
$ws.Cells.Item(70, 10).Value = 29800
$ws.Cells.Item(70, 12).Value = 29800
$ws.Cells.Item(70, 14).Value = -30340

$ws.Cells.Item(73, 8).Value = 29800
$ws.Cells.Item(73, 10).Value = 29800
$ws.Cells.Item(73, 12).Value = 29800
$ws.Cells.Item(73, 14).Value = -31672

$ws.Cells.Item(100, 8).Value = 3413.95
$ws.Cells.Item(100, 9).Value = 3031.111
$ws.Cells.Item(100, 10).Value = 3727.182
$ws.Cells.Item(100, 11).Value = 3031.111
$ws.Cells.Item(100, 12).Value = 3727.182
$ws.Cells.Item(100, 13).Value = -2490.111
$ws.Cells.Item(100, 14).Value = -4809.182

$ws.Cells.Item(132, 8).Value = 20840890
$ws.Cells.Item(132, 9).Value = 43480144
$ws.Cells.Item(132, 10).Value = 12776.4
$ws.Cells.Item(132, 11).Value = 130440432
$ws.Cells.Item(132, 12).Value = 38329.2
$ws.Cells.Item(132, 13).Value = -130437902
$ws.Cells.Item(132, 14).Value = -43389.2

$ws.Cells.Item(136, 8).Value = 3055.5
$ws.Cells.Item(136, 10).Value = 3846.6924
$ws.Cells.Item(136, 12).Value = 11540.0772
$ws.Cells.Item(136, 14).Value = -16640.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(70, 8).Value = 22808.46
$ws.Cells.Item(70, 10).Value = 24200.834
$ws.Cells.Item(70, 12).Value = 24200.834
$ws.Cells.Item(70, 14).Value = -24830.834

$ws.Cells.Item(73, 8).Value = 22808.46
$ws.Cells.Item(73, 10).Value = 24200.834
$ws.Cells.Item(73, 12).Value = 24200.834
$ws.Cells.Item(73, 14).Value = -26384.834

$ws.Cells.Item(107, 8).Value = 783.7406999999999
$ws.Cells.Item(107, 9).Value = 1271.7
$ws.Cells.Item(107, 10).Value = 496.70587
$ws.Cells.Item(107, 11).Value = 3815.1
$ws.Cells.Item(107, 12).Value = 1490.11761
$ws.Cells.Item(107, 13).Value = -1895.1
$ws.Cells.Item(107, 14).Value = -5330.11761

$ws.Cells.Item(122, 8).Value = 1869.6842
$ws.Cells.Item(122, 9).Value = 1333.5652
$ws.Cells.Item(122, 11).Value = 4000.6956
$ws.Cells.Item(122, 13).Value = -1550.6956

$ws.Cells.Item(136, 8).Value = 1953.7935
$ws.Cells.Item(136, 9).Value = 2248.5833
$ws.Cells.Item(136, 10).Value = 1401.0625
$ws.Cells.Item(136, 11).Value = 6745.749899999999
$ws.Cells.Item(136, 12).Value = 4203.1875
$ws.Cells.Item(136, 13).Value = -4195.749899999999
$ws.Cells.Item(136, 14).Value = -9303.1875

Write-Output "All edits applied"